$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# NumberFormat is set to Text ("@") before assigning the value so that
# Excel keeps the original literal text (e.g. "1.44%") instead of
# auto-converting it into a numeric percentage value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.44%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.78%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.122"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.30%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07853"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.98%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.251"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-6.68%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.814"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.04%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.809"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.13%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.76%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1777"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.58%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07674"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.26%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08893"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.27%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03102"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.36%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.55%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001505"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.24%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005910"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.94%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.467"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.29%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.250"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.11%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.23%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.75%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.316"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-6.72%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "10.67%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04607"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.40%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.33%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004482"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.07%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001249"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.85%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-1.42%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.26%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04792"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "7.18%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007290"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.71%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1364"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.54%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002188"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.54%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.10%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006278"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.89%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.27%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.002509"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-52.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.065"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "29.83%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.27%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.27%"
